$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("week 13")

# Insert 9 new rows above the old totals row (row 8), pushing it down to row 17
$ws.Rows.Item(8).Resize(9).Insert()

# Copy formatting from the existing data row (row 7) down onto the freshly inserted rows
$ws.Range("A7:F7").Copy()
$ws.Range("A8:F16").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# B column date-text cells use a different display number format than the time columns;
# match the target layout (m/d/yyyy for most, m/d/yyyy h:mm for the very last row)
$ws.Range("B8:B15").NumberFormat = "m/d/yyyy"
$ws.Range("B16").NumberFormat = "m/d/yyyy h:mm"

Write-Output "step1 done"
